# Update PFAS summary values on the "sum_levels_ng_gdw" and
# "sum_levels_ng_gww" sheets (rows 10-17, columns C/D/E).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("sum_levels_ng_gdw")
$ws2 = $wb.Worksheets.Item("sum_levels_ng_gww")

# --- sum_levels_ng_gdw ---------------------------------------------------
$ws1.Range("C10").Value = 18.45
$ws1.Range("D10").Value = 6.73
$ws1.Range("E10").Value = 54.8

$ws1.Range("C11").Value = 10.67
$ws1.Range("D11").Value = 2.93
$ws1.Range("E11").Value = 144.48

$ws1.Range("C12").Value = 6.32
$ws1.Range("D12").Value = 3.5
$ws1.Range("E12").Value = 15.9

$ws1.Range("C13").Value = 15.87
$ws1.Range("D13").Value = 9.369999999999999
$ws1.Range("E13").Value = 25.33

$ws1.Range("C14").Value = 23.86
$ws1.Range("D14").Value = 2.93
$ws1.Range("E14").Value = 144.48

$ws1.Range("C15").Value = 5.66
$ws1.Range("D15").Value = 2.93
$ws1.Range("E15").Value = 15.9

$ws1.Range("C16").Value = 10.36
$ws1.Range("D16").Value = 3.5
$ws1.Range("E16").Value = 28.97

$ws1.Range("C17").Value = 18.45
$ws1.Range("D17").Value = 6.73
$ws1.Range("E17").Value = 144.48

# --- sum_levels_ng_gww ---------------------------------------------------
$ws2.Range("C10").Value = 4.13
$ws2.Range("D10").Value = 1.55
$ws2.Range("E10").Value = 11.84

$ws2.Range("C11").Value = 1.96
$ws2.Range("D11").Value = 0.42
$ws2.Range("E11").Value = 19.22

$ws2.Range("C12").Value = 0.98
$ws2.Range("D12").Value = 0.42
$ws2.Range("E12").Value = 3

$ws2.Range("C13").Value = 2.34
$ws2.Range("D13").Value = 1.97
$ws2.Range("E13").Value = 4.97

$ws2.Range("C14").Value = 4.92
$ws2.Range("D14").Value = 0.64
$ws2.Range("E14").Value = 19.22

$ws2.Range("C15").Value = 0.9399999999999999
$ws2.Range("D15").Value = 0.42
$ws2.Range("E15").Value = 3

$ws2.Range("C16").Value = 1.38
$ws2.Range("D16").Value = 0.5600000000000001
$ws2.Range("E16").Value = 6.23

$ws2.Range("C17").Value = 4.36
$ws2.Range("D17").Value = 1.55
$ws2.Range("E17").Value = 19.22
